$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column J -- copy formatting from the neighboring header
# cell (I1) so the new header gets the same bold/centered/bordered style.
$ws.Range("J1").Value = "Q8"
$ws.Range("I1").Copy()
$ws.Range("J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New values for row 4 (G4:J4)
$ws.Range("G4").Value = 0.2890697267702507
$ws.Range("H4").Value = -0.6507920071323952
$ws.Range("I4").Value = 0.4578003130087183
$ws.Range("J4").Value = -0.1119550751434417

# New values for row 8 (G8:I8)
$ws.Range("G8").Value = 0.7010458975705092
$ws.Range("H8").Value = 0.6218889942996384
$ws.Range("I8").Value = 0.4230596606995932
